$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

Replace-Text "Estamos a remover o Tether Omni (USDT) a 29 de setembro" "Vamos remover o Tether Omni (USDT) no dia 29 de setembro"

Replace-Text "Deixaremos de oferecer o Tether Omni (USDT) como moeda da conta no Deriv, a partir de 29/09/2023 (00:00 GMT). Isso ocorre porque o Tether parou de oferecer suporte ao Omni para transferências de USDT." "Vamos deixar de oferecer o Tether Omni (USDT) como moeda da conta no Deriv, a partir de 29/09/2023 (00:00 GMT). O motivo deve-se ao facto de a Tether ter deixado de suportar a Omni para transferências de USDT."

Replace-Text "O que é que preciso de fazer?" "O que é necessário fazer? "

Replace-Text ", levante o seu saldo antes da data acima indicada. Se tiver posições abertas, feche-as primeiro antes de levantar o seu saldo." ", retire o seu saldo antes da data acima indicada. Caso tenha posições abertas, feche-as primeiro antes de retirar o seu saldo."

Replace-Text "A sua conta USDT será encerrada em 29/09/2023 às 00:00 GMT. Todas as posições abertas serão automaticamente fechadas e o saldo da conta será transferido para a sua última conta ativa após a data mencionada" "A sua conta USDT será encerrada no dia 29/09/2023 às 00:00 GMT. Todas as posições abertas serão automaticamente fechadas e o saldo da conta será transferido para a sua última conta ativa após a data mencionada"
